$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47 and 48 swap: Stacks <-> ApeXProtocol content, plus all market data refresh

$ws.Range("D2").Value = '42.827.14'
$ws.Range("E2").Value = '  -6.86%  '
$ws.Range("D3").Value = '2.546.12'
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '''296.81'
$ws.Range("E5").Value = '  -4.42%  '
$ws.Range("D6").Value = '''93.46'
$ws.Range("E6").Value = '  -4.95%  '
$ws.Range("D7").Value = '''0.573'
$ws.Range("E7").Value = '  -4.28%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -5.59%  '
$ws.Range("D10").Value = '''35.62'
$ws.Range("E10").Value = '  -8.09%  '
$ws.Range("D11").Value = '''0.0807'
$ws.Range("E11").Value = '  -3.56%  '
$ws.Range("D12").Value = '''7.67'
$ws.Range("E12").Value = '  -5.32%  '
$ws.Range("D13").Value = '2.934.71'
$ws.Range("E13").Value = '  -2.02%  '
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("D15").Value = '2.540.45'
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("E16").Value = '  -5.55%  '
$ws.Range("D17").Value = '''14.11'
$ws.Range("E17").Value = '  -4.75%  '
$ws.Range("D18").Value = '42.829.02'
$ws.Range("E18").Value = '  -7.24%  '
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("E20").Value = '  -3.87%  '
$ws.Range("D21").Value = '''12.44'
$ws.Range("E21").Value = '  -2.35%  '
$ws.Range("D22").Value = '''72.33'
$ws.Range("E22").Value = '  -0.56%  '
$ws.Range("D23").Value = '''260.25'
$ws.Range("E23").Value = '  -11.47%  '
$ws.Range("E24").Value = '  -5.52%  '
$ws.Range("D25").Value = '''29.51'
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("E26").Value = '  -5.68%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").Value = '''10.00'
$ws.Range("E28").Value = '  -7.00%  '
$ws.Range("E29").Value = '  -4.36%  '
$ws.Range("D30").Value = '''36.05'
$ws.Range("E30").Value = '  -5.65%  '
$ws.Range("D31").Value = '''5.91'
$ws.Range("E31").Value = '  -4.86%  '
$ws.Range("D32").Value = '''150.67'
$ws.Range("E32").Value = '  -2.92%  '
$ws.Range("D33").Value = '''2.16'
$ws.Range("E33").Value = '  -1.66%  '
$ws.Range("E34").Value = '  -5.23%  '
$ws.Range("E35").Value = '  -1.86%  '
$ws.Range("E36").Value = '  -5.13%  '
$ws.Range("E37").Value = '  -6.55%  '
$ws.Range("D38").Value = '''24.18'
$ws.Range("E38").Value = '  +15.03%  '
$ws.Range("E39").Value = '  -2.96%  '
$ws.Range("D40").Value = '''16.21'
$ws.Range("E40").Value = '  +3.47%  '
$ws.Range("D41").Value = '''3.40'
$ws.Range("E41").Value = '  -4.51%  '
$ws.Range("E42").Value = '  -6.22%  '
$ws.Range("E43").Value = '  -3.17%  '
$ws.Range("D44").Value = '2.071.45'
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("D45").Value = '''0.998'
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("D46").Value = '''84.82'
$ws.Range("E46").Value = '  -12.96%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = '''1.72'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").Value = '''1.58'
$ws.Range("E48").Value = '  +2.73%  '
$ws.Range("D49").Value = '2.791.48'
$ws.Range("E49").Value = '  -2.09%  '
$ws.Range("D50").Value = '''103.76'
$ws.Range("E50").Value = '  -3.81%  '
$ws.Range("D51").Value = '''8.63'
$ws.Range("E51").Value = '  -10.05%  '
